$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 132.1614712523325
$ws.Range("C3").Value = 12.63512059316575
$ws.Range("C4").Value = 12.73197136608733
$ws.Range("C5").Value = 16.72784680372324
$ws.Range("C6").Value = 23.61518644584078
$ws.Range("C7").Value = 7.377060889067842
$ws.Range("C8").Value = 7.523899157690896
$ws.Range("C9").Value = 22.1460227049899
$ws.Range("C10").Value = 37.30316866721864
$ws.Range("C11").Value = 10.42707918147757
$ws.Range("C12").Value = 2.239283596501588
$ws.Range("C13").Value = 5.99381315645385
$ws.Range("C14").Value = 1.418395190529084
$ws.Range("C15").Value = 2.500936894313947
$ws.Range("C16").Value = 18.62580953113826
$ws.Range("C17").Value = 19.63493210061159
$ws.Range("C18").Value = 18.9819704380112
$ws.Range("C19").Value = 6.045362761395987
$ws.Range("C20").Value = 29.43716758582084
$ws.Range("C21").Value = 69.8075377470966
$ws.Range("C22").Value = 13.17404828119717
$ws.Range("C23").Value = 2.510309549757971
$ws.Range("C24").Value = 22.95910056475905
$ws.Range("C25").Value = 6.689732823172689
$ws.Range("C26").Value = 13.34353713380995
$ws.Range("C27").Value = 27.71962847570331
$ws.Range("C28").Value = 4.23019182373652
$ws.Range("C29").Value = 14.89314950055538
$ws.Range("C30").Value = 2.409553503734705
$ws.Range("C31").Value = 2.574356028625474
$ws.Range("C32").Value = 4.396556457867959
$ws.Range("C33").Value = 5.243219666311532
$ws.Range("C34").Value = 105.6509153289084
$ws.Range("C35").Value = 9.564013826006956
$ws.Range("C36").Value = 23.62299699204414
$ws.Range("C37").Value = 4.27939826481765
$ws.Range("C38").Value = 9.728816350897727
$ws.Range("C39").Value = 9.168019133496911
$ws.Range("C40").Value = 7.915988577099266
$ws.Range("C41").Value = 5.656397560468958
$ws.Range("C42").Value = 246.936034093561
